# Update "想去人数" (F column) values across the workbook's sheets.
# Sheet 1 = 展览 (Exhibitions), Sheet 3 = 本地生活 (Local Life),
# Sheet 4 = 全部类型 (All Types, an aggregate of all category sheets).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$ws1.Range("F5").Value  = 763
$ws1.Range("F6").Value  = 2463
$ws1.Range("F7").Value  = 54
$ws1.Range("F8").Value  = 1829
$ws1.Range("F9").Value  = 3104
$ws1.Range("F10").Value = 190
$ws1.Range("F11").Value = 4600
$ws1.Range("F12").Value = 421
$ws1.Range("F18").Value = 260
$ws1.Range("F22").Value = 319
$ws1.Range("F23").Value = 4599
$ws1.Range("F25").Value = 23
$ws1.Range("F27").Value = 4730
$ws1.Range("F29").Value = 1158
$ws1.Range("F31").Value = 618
$ws1.Range("F33").Value = 49
$ws1.Range("F35").Value = 712
$ws1.Range("F37").Value = 657
$ws1.Range("F38").Value = 650

# --- Sheet 3: 本地生活 ---
$ws3.Range("F3").Value = 1057

# --- Sheet 4: 全部类型 ---
$ws4.Range("F4").Value  = 1057
$ws4.Range("F8").Value  = 763
$ws4.Range("F9").Value  = 2463
$ws4.Range("F10").Value = 54
$ws4.Range("F11").Value = 1829
$ws4.Range("F13").Value = 3104
$ws4.Range("F14").Value = 190
$ws4.Range("F15").Value = 4600
$ws4.Range("F16").Value = 421
$ws4.Range("F22").Value = 260
$ws4.Range("F27").Value = 319
$ws4.Range("F28").Value = 4599
$ws4.Range("F30").Value = 23
$ws4.Range("F32").Value = 4731
$ws4.Range("F34").Value = 1158
$ws4.Range("F36").Value = 618
$ws4.Range("F39").Value = 49
$ws4.Range("F41").Value = 712
$ws4.Range("F43").Value = 657
$ws4.Range("F44").Value = 650
